$d = $word.ActiveDocument

# Locate the question paragraph "2.5.6 ..." and switch it to the
# "Заголовок-2" paragraph style (replacing its direct spacing/indent
# formatting with the style reference, as Word does when a style is applied).
$rng = $d.Content
$rng.Find.Execute("2.5.6 Які розділи містить технічне завдання?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$questionPara = $rng.Paragraphs(1)
$questionPara.Style = "Заголовок-2"

# The next (currently empty) paragraph gets the new answer text. Seed it
# with a duplicate of the question run's formatting (same character style +
# shading) via FormattedText, then overwrite the text itself so the new
# run keeps that exact run formatting.
$answerPara = $questionPara.Next()
$answerRange = $answerPara.Range
$answerRange.FormattedText = $questionPara.Range.Duplicate.FormattedText
$answerRange.Text = "Технічне завдання містить кілька розділів, а саме: вступ, підстави для розробки, призначення розробки, вимоги до програми чи програмного виробу, вимоги до програмної документації, техніко-економічні показники, стадії та етапи розробки, порядок контролю та приймання."
